$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "249.71"
Set-TextValue 2 7 "19"
Set-TextValue 3 4 "22.98"
Set-TextValue 3 7 "19"
Set-TextValue 4 4 "5.422"
Set-TextValue 4 7 "19"
Set-TextValue 5 4 "0.05638"
Set-TextValue 5 7 "19"
Set-TextValue 6 4 "3.426"
Set-TextValue 6 7 "19"
Set-TextValue 7 4 "6.369"
Set-TextValue 7 7 "19"
Set-TextValue 8 4 "0.8153"
Set-TextValue 8 7 "19"
Set-TextValue 9 4 "0.9210"
Set-TextValue 9 7 "19"
Set-TextValue 10 2 "One"
Set-TextValue 10 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 10 4 "0.01154"
Set-TextValue 10 5 "9OneONE"
Set-TextValue 10 7 "19"
Set-TextValue 11 2 "WazirX"
Set-TextValue 11 3 "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue 11 4 "0.1438"
Set-TextValue 11 5 "10WazirXWRX"
Set-TextValue 11 7 "19"
Set-TextValue 12 2 "MandalaExchangeToken"
Set-TextValue 12 3 "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue 12 4 "0.07511"
Set-TextValue 12 5 "11MandalaExchangeTokenMDX"
Set-TextValue 12 7 "19"
Set-TextValue 13 2 "LiechtensteinCryptoassetsExchange"
Set-TextValue 13 3 "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue 13 4 "0.03131"
Set-TextValue 13 5 "12LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue 13 7 "19"
Set-TextValue 14 2 "BitrueCoin"
Set-TextValue 14 3 "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue 14 4 "0.03109"
Set-TextValue 14 5 "13BitrueCoinBTR"
Set-TextValue 14 7 "19"
Set-TextValue 15 2 "BitMartToken"
Set-TextValue 15 3 "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue 15 4 "0.09341"
Set-TextValue 15 5 "14BitMartTokenBMX"
Set-TextValue 15 7 "19"
Set-TextValue 16 2 "MCDex"
Set-TextValue 16 3 "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue 16 4 "3.553"
Set-TextValue 16 5 "15MCDexMCB"
Set-TextValue 16 7 "19"
Set-TextValue 17 2 "BitForexToken"
Set-TextValue 17 3 "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue 17 4 "0.001584"
Set-TextValue 17 5 "16BitForexTokenBF"
Set-TextValue 17 7 "19"
Set-TextValue 18 2 "CoinExToken"
Set-TextValue 18 3 "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue 18 4 "0.04756"
Set-TextValue 18 5 "17CoinExTokenCET"
Set-TextValue 18 7 "19"
Set-TextValue 19 4 "0.006381"
Set-TextValue 19 7 "19"
Set-TextValue 20 4 "0.005006"
Set-TextValue 20 7 "19"
Set-TextValue 21 7 "19"
Set-TextValue 22 4 "0.0001499"
Set-TextValue 22 7 "19"
Set-TextValue 23 4 "3.697"
Set-TextValue 23 7 "19"
Set-TextValue 24 4 "2.179"
Set-TextValue 24 7 "19"
Set-TextValue 25 4 "0.3298"
Set-TextValue 25 7 "19"
Set-TextValue 26 7 "19"
Set-TextValue 27 7 "19"
Set-TextValue 28 4 "0.0003030"
Set-TextValue 28 7 "19"
Set-TextValue 29 7 "19"
Set-TextValue 30 7 "19"
Set-TextValue 31 7 "19"
Set-TextValue 32 7 "19"
Set-TextValue 33 7 "19"
Set-TextValue 34 7 "19"
Set-TextValue 35 7 "19"
Set-TextValue 36 7 "19"
Set-TextValue 37 7 "19"
Set-TextValue 38 7 "19"
Set-TextValue 39 7 "19"
Set-TextValue 40 4 "0.04047"
Set-TextValue 40 7 "19"
Set-TextValue 41 4 "0.006796"
Set-TextValue 41 7 "19"
Set-TextValue 42 4 "0.1070"
Set-TextValue 42 7 "19"
Set-TextValue 43 4 "0.002719"
Set-TextValue 43 7 "19"
Set-TextValue 44 4 "0.007556"
Set-TextValue 44 7 "19"
Set-TextValue 45 4 "0.00005800"
Set-TextValue 45 7 "19"
Set-TextValue 46 7 "19"
Set-TextValue 47 4 "0.4997"
Set-TextValue 47 7 "19"
Set-TextValue 48 4 "0.2411"
Set-TextValue 48 7 "19"
Set-TextValue 49 7 "19"
Set-TextValue 50 4 "0.01009"
Set-TextValue 50 7 "19"
Set-TextValue 51 7 "19"
